$d = $word.ActiveDocument

$r = $d.Paragraphs.Item(3).Range
$r.Font.Size = 10
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.StrikeThrough = $False

$r = $d.Paragraphs.Item(4).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.StrikeThrough = $False

$r = $d.Paragraphs.Item(5).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False

$r = $d.Paragraphs.Item(6).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.StrikeThrough = $True

$r = $d.Paragraphs.Item(8).Range
$r.Font.Size = 10
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.StrikeThrough = $False

$r = $d.Paragraphs.Item(9).Range
$r.Font.Size = 5
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.StrikeThrough = $False
$r.Font.Color = 8323327

$r = $d.Paragraphs.Item(10).Range
$r.Font.Size = 6
$r.Font.Bold = $True
$r.Font.Italic = $True
$r.Font.StrikeThrough = $False
$r.Font.Color = 16744192

$r = $d.Paragraphs.Item(11).Range
$r.Font.Size = 7
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False
$r.Font.Color = 65407

$r = $d.Paragraphs.Item(12).Range
$r.Font.Size = 8
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.StrikeThrough = $True
$r.Font.Color = 8323327

$r = $d.Paragraphs.Item(14).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.StrikeThrough = $False

$r = $d.Paragraphs.Item(15).Range
$r.Font.Size = 9
$r.Font.Bold = $True
$r.Font.Italic = $True
$r.Font.StrikeThrough = $False
$r.Font.Color = 16744192

$r = $d.Paragraphs.Item(16).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.StrikeThrough = $False
$r.Font.Color = 65407

$r = $d.Paragraphs.Item(17).Range
$r.Font.Size = 11
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False
$r.Font.Color = 8323327

$r = $d.Paragraphs.Item(18).Range
$r.Font.Size = 12
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.StrikeThrough = $True
$r.Font.Color = 16744192

$r = $d.Paragraphs.Item(20).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False

$r = $d.Paragraphs.Item(21).Range
$r.Font.Size = 13
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False
$r.Font.Color = 65407

$r = $d.Paragraphs.Item(22).Range
$r.Font.Size = 14
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False
$r.Font.Color = 8323327

$r = $d.Paragraphs.Item(23).Range
$r.Font.Size = 15
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $False
$r.Font.Color = 16744192

$r = $d.Paragraphs.Item(24).Range
$r.Font.Size = 16
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $True
$r.Font.Color = 65407

$r = $d.Paragraphs.Item(26).Range
$r.Font.Size = 10
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.StrikeThrough = $True

$r = $d.Paragraphs.Item(27).Range
$r.Font.Size = 17
$r.Font.Bold = $True
$r.Font.Italic = $False
$r.Font.StrikeThrough = $True
$r.Font.Color = 8323327

$r = $d.Paragraphs.Item(28).Range
$r.Font.Size = 18
$r.Font.Bold = $False
$r.Font.Italic = $True
$r.Font.StrikeThrough = $True
$r.Font.Color = 16744192

$r = $d.Paragraphs.Item(29).Range
$r.Font.Size = 19
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.Underline = 1
$r.Font.StrikeThrough = $True
$r.Font.Color = 65407

$r = $d.Paragraphs.Item(30).Range
$r.Font.Size = 20
$r.Font.Bold = $False
$r.Font.Italic = $False
$r.Font.StrikeThrough = $True
$r.Font.Color = 8323327
